# Add 2019 data (column AD) for OOSS y huelgas cuadro7 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, matching the style of the other year headers (e.g. AC1):
# bold font, centered horizontal alignment (xlCenter = -4108)
$ws.Range("AD1").Value = "Tamaño_2019"
$ws.Range("AD1").Font.Bold = $true
$ws.Range("AD1").HorizontalAlignment = -4108

# New data values for 2019
$ws.Range("AD2").Value = 99.3
$ws.Range("AD3").Value = 50.7
$ws.Range("AD4").Value = 120
$ws.Range("AD5").Value = 124.7
$ws.Range("AD6").Value = 94.8
$ws.Range("AD7").Value = 277
$ws.Range("AD8").Value = 90.2
$ws.Range("AD9").Value = 204.6
$ws.Range("AD10").Value = 102.1
$ws.Range("AD11").Value = 100
$ws.Range("AD12").Value = 77.59999999999999
